$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NIFTY")

# Row 7 - summary/header values
$ws.Range("F7").Value = 21641.8
$ws.Range("G7").Value = 21680
$ws.Range("H7").Value = 21164.15
$ws.Range("I7").Value = 21217.05
$ws.Range("J7").Value = 21217.05

# Row 9
$ws.Range("G9").Value = 21158
$ws.Range("H9").Value = 21060
$ws.Range("I9").Value = 21101.6

# Row 10
$ws.Range("G10").Value = 21173.9
$ws.Range("H10").Value = 21081.1
$ws.Range("I10").Value = 21157.25

# Row 11
$ws.Range("G11").Value = 21228.3
$ws.Range("H11").Value = 21141.1
$ws.Range("I11").Value = 21190.05

# Row 12
$ws.Range("G12").Value = 21244.95
$ws.Range("H12").Value = 21140.1
$ws.Range("I12").Value = 21236

# Row 13
$ws.Range("G13").Value = 21258.55
$ws.Range("H13").Value = 21195.7
$ws.Range("I13").Value = 21218

# Row 14
$ws.Range("G14").Value = 21250.1
$ws.Range("H14").Value = 21215
$ws.Range("I14").Value = 21240

# Row 15
$ws.Range("G15").Value = 21289.5
$ws.Range("H15").Value = 21235
$ws.Range("I15").Value = 21285.4

# Row 16
$ws.Range("G16").Value = 21340
$ws.Range("H16").Value = 21281.6
$ws.Range("I16").Value = 21325.25

# Row 17
$ws.Range("G17").Value = 21333
$ws.Range("H17").Value = 21297.05
$ws.Range("I17").Value = 21321.1

# Row 18
$ws.Range("G18").Value = 21348
$ws.Range("H18").Value = 21255.6
$ws.Range("I18").Value = 21294.35

# Row 19
$ws.Range("G19").Value = 21328
$ws.Range("H19").Value = 21258
$ws.Range("I19").Value = 21324.15

# Row 20
$ws.Range("G20").Value = 21375
$ws.Range("H20").Value = 21317.7
$ws.Range("I20").Value = 21345.9

# Row 21
$ws.Range("G21").Value = 21379.7
$ws.Range("H21").Value = 21313.1
$ws.Range("I21").Value = 21376
